$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-10
# from serial date 45208 (2023-10-09) to 45212 (2023-10-13)
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
